$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting text for the "R10" row (cell E8) to reflect the commit update.
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the active selection change recorded in the sheet view.
$ws.Range("E8").Select()
